$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the conversion text in A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 7.46 = 30084.04 pesos`n✅ 30084.04 pesos = 7.43 = 922.54 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$wsHoja1.Range("A1").Value = $newText

# --- Sheet "tasas": update the N10/O10/N12/O12 rate values ---
$wsTasas = $wb.Worksheets.Item("tasas")

$wsTasas.Range("N10").Value = 133.99
$wsTasas.Range("O10").Value = 4030.96
$wsTasas.Range("N12").Value = 4047.15
$wsTasas.Range("O12").Value = 124.107
